# Apply updates to parameters_0.xlsx as described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Cell value (parameter bound) changes ---
# Row 5 (dCit): bmax 0.1 -> 1E-4
$ws.Range("C5").Value = 0.0001

# Row 6 (indTime): bmin 0.01 -> 1E-3
$ws.Range("B6").Value = 0.001

# Row 7 (mu): bmin 1 -> 0.1
$ws.Range("B7").Value = 0.1

# Row 8 (nMperUnit): bmax 20 -> 10
$ws.Range("C8").Value = 10

# Row 12 (KdLacI): bmin 1E-3 -> 1E-4
$ws.Range("B12").Value = 0.0001

# Row 19 (k_{LacI}): bmin 0.1 -> 0.01, bmax 100 -> 10
$ws.Range("B19").Value = 0.01
$ws.Range("C19").Value = 10

# Row 22 (k_{LacI_W220F_Q60G_T167A}): bmin 1E-3 -> 1E-4, bmax 10 -> 1
$ws.Range("B22").Value = 0.0001
$ws.Range("C22").Value = 1

# --- View / selection changes on the worksheet ---
# Previously the sheet view was scrolled so row 3 was at the top and A25 was
# selected; now the sheet is scrolled back to the top and C22 (the last
# parameter bound that was edited) is the selected cell.
$ws.Activate()
$ws.Range("C22").Select()
